# Insert a new "Industry" column as column C, shifting the existing
# Mutual Fund / Status / Jan_2026 / Dec_2025 / Oct_2025 / MoM / QoQ
# columns one place to the right (C:I -> D:J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing columns C..I right by inserting a new blank column at C.
$ws.Range("C:C").Insert()

# Header
$ws.Range("C1").Value = "Industry"

# Per-row industry classification
$industry = @{
    2  = "Insurance"
    3  = "Auto Components"
    4  = "Banks"
    5  = "Banks"
    6  = "Consumable Fuels"
    7  = "Banks"
    8  = "Pharmaceuticals & Biotechnology"
    9  = "Leisure Services"
    10 = "Automobiles"
    11 = "Power"
    12 = "Industrial Products"
    13 = "Construction"
    14 = "Power"
    15 = "Pharmaceuticals & Biotechnology"
    16 = "Finance"
    17 = "Banks"
    18 = "Realty"
    19 = "Diversified FMCG"
    20 = "Banks"
    21 = "Realty"
}

foreach ($row in $industry.Keys) {
    $ws.Cells.Item($row, 3).Value = $industry[$row]
}
